$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.158.71'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.789.87'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +21.70%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '616.92'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +7.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.14'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.783.77'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +21.51%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.549'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +6.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +10.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.40'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.501'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +7.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.59'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +11.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000258'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +6.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.435.17'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +22.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.803.59'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +22.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.353.54'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.54%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.61'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +8.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '525.17'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.90%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.71'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.49'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +23.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.745'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +8.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.61'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.87%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.49'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +9.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.55'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +6.23%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +4.50%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000123'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +31.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.52'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +8.65%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.88'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +10.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.92'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.15'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +14.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.115'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.57%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +11.41%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +10.61%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +7.85%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +8.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.16'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +7.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.70'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.87'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.73%  '
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.146.46'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +12.96%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '427.96'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +14.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '44.53'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -7.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.75'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0368'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +6.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.84'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +4.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.53'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +7.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.95'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.99%  '
